# Applies the "Updated cryptos list" refresh: new Price (col D) and
# Volume(1h) (col E) readings for each coin row on Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'" + '64.221.56'
$ws.Range("E2").Value = '  -3.67%  '
$ws.Range("D3").Value = "'" + '3.149.81'
$ws.Range("E3").Value = '  -9.13%  '
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").Value = "'" + '562.41'
$ws.Range("E5").Value = '  -3.93%  '
$ws.Range("D6").Value = "'" + '167.91'
$ws.Range("E6").Value = '  -5.65%  '
$ws.Range("D7").Value = "'" + '0.615'
$ws.Range("E7").Value = '  -2.16%  '
$ws.Range("E8").Value = '  +0.01%  '
$ws.Range("D9").Value = "'" + '3.147.73'
$ws.Range("E9").Value = '  -9.07%  '
$ws.Range("D10").Value = "'" + '0.124'
$ws.Range("E10").Value = '  -7.38%  '
$ws.Range("E11").Value = '  -6.63%  '
$ws.Range("E12").Value = '  -6.58%  '
$ws.Range("D13").Value = "'" + '3.692.41'
$ws.Range("E13").Value = '  -9.03%  '
$ws.Range("E14").Value = '  +0.72%  '
$ws.Range("D15").Value = "'" + '27.00'
$ws.Range("E15").Value = '  -10.39%  '
$ws.Range("D16").Value = "'" + '64.171.98'
$ws.Range("E16").Value = '  -3.54%  '
$ws.Range("D17").Value = "'" + '0.0000163'
$ws.Range("E17").Value = '  -6.57%  '
$ws.Range("D18").Value = "'" + '3.143.43'
$ws.Range("E18").Value = '  -8.86%  '
$ws.Range("D19").Value = "'" + '5.71'
$ws.Range("E19").Value = '  -4.71%  '
$ws.Range("D20").Value = "'" + '12.85'
$ws.Range("E20").Value = '  -7.71%  '
$ws.Range("D21").Value = "'" + '350.32'
$ws.Range("E21").Value = '  -5.88%  '
$ws.Range("D22").Value = "'" + '7.17'
$ws.Range("E22").Value = '  -6.82%  '
$ws.Range("D23").Value = "'" + '0.999'
$ws.Range("E23").Value = '  -0.10%  '
$ws.Range("D24").Value = "'" + '67.92'
$ws.Range("E24").Value = '  -7.66%  '
$ws.Range("E25").Value = '  -7.59%  '
$ws.Range("D26").Value = "'" + '0.0000115'
$ws.Range("E26").Value = '  -9.37%  '
$ws.Range("D27").Value = "'" + '9.56'
$ws.Range("E27").Value = '  -4.53%  '
$ws.Range("E28").Value = '  -1.61%  '
$ws.Range("D29").Value = "'" + '0.999'
$ws.Range("E29").Value = '  -0.08%  '
$ws.Range("E30").Value = '  -0.16%  '
$ws.Range("E31").Value = '  -6.25%  '
$ws.Range("D32").Value = "'" + '5.43'
$ws.Range("E32").Value = '  -9.31%  '
$ws.Range("D33").Value = "'" + '21.85'
$ws.Range("E33").Value = '  -7.86%  '
$ws.Range("D34").Value = "'" + '6.57'
$ws.Range("E34").Value = '  -7.40%  '
$ws.Range("D35").Value = "'" + '1.19'
$ws.Range("E35").Value = '  -7.14%  '
$ws.Range("D36").Value = "'" + '153.67'
$ws.Range("E36").Value = '  -5.39%  '
$ws.Range("D37").Value = "'" + '1.42'
$ws.Range("E37").Value = '  -10.24%  '
$ws.Range("D38").Value = "'" + '0.814'
$ws.Range("E38").Value = '  -8.21%  '
$ws.Range("D39").Value = "'" + '26.13'
$ws.Range("E39").Value = '  -6.55%  '
$ws.Range("D40").Value = "'" + '1.69'
$ws.Range("E40").Value = '  -7.17%  '
$ws.Range("D41").Value = "'" + '2.43'
$ws.Range("E41").Value = '  -5.72%  '
$ws.Range("D42").Value = "'" + '2.598.75'
$ws.Range("E42").Value = '  -6.38%  '
$ws.Range("D43").Value = "'" + '4.15'
$ws.Range("E43").Value = '  -8.35%  '
$ws.Range("E44").Value = '  -2.14%  '
$ws.Range("E45").Value = '  -8.40%  '
$ws.Range("D46").Value = "'" + '0.0644'
$ws.Range("E46").Value = '  -7.53%  '
$ws.Range("D47").Value = "'" + '23.69'
$ws.Range("E47").Value = '  -7.16%  '
$ws.Range("D48").Value = "'" + '316.16'
$ws.Range("E48").Value = '  -7.40%  '
$ws.Range("D49").Value = "'" + '0.0270'
$ws.Range("E49").Value = '  -6.73%  '
$ws.Range("D50").Value = "'" + '0.102'
$ws.Range("E50").Value = '  -3.29%  '
$ws.Range("E51").Value = '  -0.03%  '

# Entering text with a leading apostrophe makes Excel mark the cells with a
# "number stored as text" quote-prefix format; reset to Normal so the cell
# styling matches the original (unstyled) cells.
$ws.Range("D2:D51").Style = "Normal"

